$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the vm_pu values per the diff (rows 2-25, columns B-N excluding G/H)
# Row 2
$ws.Range("B2").Value2 = 1.02
$ws.Range("C2").Value2 = 1.02523330821534
$ws.Range("D2").Value2 = 1.049233876070395
$ws.Range("E2").Value2 = 1.025619745598501
$ws.Range("F2").Value2 = 1.054405761363375
$ws.Range("I2").Value2 = 1.041920402842431
$ws.Range("J2").Value2 = 1.030403474244333
$ws.Range("K2").Value2 = 1.051991503114629
$ws.Range("L2").Value2 = 1.028444785040044
$ws.Range("M2").Value2 = 1.057149060950509
$ws.Range("N2").Value2 = 1.014179181493951

# Row 3
$ws.Range("B3").Value2 = 1.02
$ws.Range("C3").Value2 = 1.026010584597794
$ws.Range("D3").Value2 = 1.049819176042703
$ws.Range("E3").Value2 = 1.026275040392894
$ws.Range("F3").Value2 = 1.055132161391283
$ws.Range("I3").Value2 = 1.042108586461252
$ws.Range("J3").Value2 = 1.03082088615185
$ws.Range("K3").Value2 = 1.052389594658597
$ws.Range("L3").Value2 = 1.028907833965634
$ws.Range("M3").Value2 = 1.057688919792122
$ws.Range("N3").Value2 = 1.014318045471715

# Row 4
$ws.Range("B4").Value2 = 1.02
$ws.Range("C4").Value2 = 1.02651421367703
$ws.Range("D4").Value2 = 1.050198256978717
$ws.Range("E4").Value2 = 1.02670004095768
$ws.Range("F4").Value2 = 1.055602946767596
$ws.Range("I4").Value2 = 1.042229299698878
$ws.Range("J4").Value2 = 1.031090985947046
$ws.Range("K4").Value2 = 1.052646838168191
$ws.Range("L4").Value2 = 1.029207754063366
$ws.Range("M4").Value2 = 1.058038327910411
$ws.Range("N4").Value2 = 1.014407876906846

# Row 5
$ws.Range("B5").Value2 = 1.02
$ws.Range("C5").Value2 = 1.026726100109428
$ws.Range("D5").Value2 = 1.050357704985643
$ws.Range("E5").Value2 = 1.026878944102127
$ws.Range("F5").Value2 = 1.055801043674606
$ws.Range("I5").Value2 = 1.042279794231531
$ws.Range("J5").Value2 = 1.031204536114308
$ws.Range("K5").Value2 = 1.052754898413394
$ws.Range("L5").Value2 = 1.029333910006227
$ws.Range("M5").Value2 = 1.058185237200805
$ws.Range("N5").Value2 = 1.014445636103145

# Row 6
$ws.Range("B6").Value2 = 1.02
$ws.Range("C6").Value2 = 1.026761686147362
$ws.Range("D6").Value2 = 1.05038448177852
$ws.Range("E6").Value2 = 1.026908996328308
$ws.Range("F6").Value2 = 1.055834315432572
$ws.Range("I6").Value2 = 1.042288257603844
$ws.Range("J6").Value2 = 1.031223601661326
$ws.Range("K6").Value2 = 1.05277303717931
$ws.Range("L6").Value2 = 1.029355096179616
$ws.Range("M6").Value2 = 1.058209904926063
$ws.Range("N6").Value2 = 1.014451975676403

# Row 7
$ws.Range("B7").Value2 = 1.02
$ws.Range("C7").Value2 = 1.026517044284526
$ws.Range("D7").Value2 = 1.050200387209685
$ws.Range("E7").Value2 = 1.026702430557544
$ws.Range("F7").Value2 = 1.055605593048692
$ws.Range("I7").Value2 = 1.042229975405826
$ws.Range("J7").Value2 = 1.031092503210751
$ws.Range("K7").Value2 = 1.052648282410108
$ws.Range("L7").Value2 = 1.029209439493845
$ws.Range("M7").Value2 = 1.058040290850662
$ws.Range("N7").Value2 = 1.014408381471002

# Row 8
$ws.Range("B8").Value2 = 1.02
$ws.Range("C8").Value2 = 1.025495850642077
$ws.Range("D8").Value2 = 1.049431606667321
$ws.Range("E8").Value2 = 1.025841001245458
$ws.Range("F8").Value2 = 1.054651093770568
$ws.Range("I8").Value2 = 1.041984218075941
$ws.Range("J8").Value2 = 1.030544538442227
$ws.Range("K8").Value2 = 1.052126111220558
$ws.Range("L8").Value2 = 1.028601212336681
$ws.Range("M8").Value2 = 1.057331490412541
$ws.Range("N8").Value2 = 1.014226115672942

# Row 9
$ws.Range("B9").Value2 = 1.02
$ws.Range("C9").Value2 = 1.02370166157443
$ws.Range("D9").Value2 = 1.048079709468397
$ws.Range("E9").Value2 = 1.024330652301043
$ws.Range("F9").Value2 = 1.052975033444306
$ws.Range("I9").Value2 = 1.041543130664632
$ws.Range("J9").Value2 = 1.029579061453318
$ws.Range("K9").Value2 = 1.051203378748355
$ws.Range("L9").Value2 = 1.027531772879144
$ws.Range("M9").Value2 = 1.056083212974213
$ws.Range("N9").Value2 = 1.013904785864855

# Row 10
$ws.Range("B10").Value2 = 1.02
$ws.Range("C10").Value2 = 1.022509199422012
$ws.Range("D10").Value2 = 1.0471804450863
$ws.Range("E10").Value2 = 1.023328979076868
$ws.Range("F10").Value2 = 1.051861753992025
$ws.Range("I10").Value2 = 1.041243729016145
$ws.Range("J10").Value2 = 1.02893556050365
$ws.Range("K10").Value2 = 1.05058656715347
$ws.Range("L10").Value2 = 1.026820469105103
$ws.Range("M10").Value2 = 1.055251621457484
$ws.Range("N10").Value2 = 1.01369049039534

# Row 11
$ws.Range("B11").Value2 = 1.02
$ws.Range("C11").Value2 = 1.02199374259251
$ws.Range("D11").Value2 = 1.046791554229574
$ws.Range("E11").Value2 = 1.022896507789351
$ws.Range("F11").Value2 = 1.051380691260786
$ws.Range("I11").Value2 = 1.041112829382625
$ws.Range("J11").Value2 = 1.028656969758631
$ws.Range("K11").Value2 = 1.050319107382211
$ws.Range("L11").Value2 = 1.026512877742201
$ws.Range("M11").Value2 = 1.054891697142015
$ws.Range("N11").Value2 = 1.01359768610007

# Row 12
$ws.Range("B12").Value2 = 1.02
$ws.Range("C12").Value2 = 1.021802414154428
$ws.Range("D12").Value2 = 1.046647179574003
$ws.Range("E12").Value2 = 1.022736059991687
$ws.Range("F12").Value2 = 1.051202154894442
$ws.Range("I12").Value2 = 1.041064019572087
$ws.Range("J12").Value2 = 1.028553497518173
$ws.Range("K12").Value2 = 1.05021970594963
$ws.Range("L12").Value2 = 1.02639868737608
$ws.Range("M12").Value2 = 1.054758031193351
$ws.Range("N12").Value2 = 1.01356321298871

# Row 13
$ws.Range("B13").Value2 = 1.02
$ws.Range("C13").Value2 = 1.021843448606304
$ws.Range("D13").Value2 = 1.046678144921775
$ws.Range("E13").Value2 = 1.022770467900091
$ws.Range("F13").Value2 = 1.051240444646786
$ws.Range("I13").Value2 = 1.041074497927815
$ws.Range("J13").Value2 = 1.028575692246748
$ws.Range("K13").Value2 = 1.050241030356491
$ws.Range("L13").Value2 = 1.026423178733339
$ws.Range("M13").Value2 = 1.054786701787759
$ws.Range("N13").Value2 = 1.01357060764666

# Row 14
$ws.Range("B14").Value2 = 1.02
$ws.Range("C14").Value2 = 1.021977924555558
$ws.Range("D14").Value2 = 1.046779618598109
$ws.Range("E14").Value2 = 1.022883241210337
$ws.Range("F14").Value2 = 1.051365930279961
$ws.Range("I14").Value2 = 1.041108798575199
$ws.Range("J14").Value2 = 1.028648416524362
$ws.Range("K14").Value2 = 1.050310891941001
$ws.Range("L14").Value2 = 1.026503437442023
$ws.Range("M14").Value2 = 1.054880647732381
$ws.Range("N14").Value2 = 1.01359483656902

# Row 15
$ws.Range("B15").Value2 = 1.02
$ws.Range("C15").Value2 = 1.022060797595399
$ws.Range("D15").Value2 = 1.046842150109207
$ws.Range("E15").Value2 = 1.022952749977446
$ws.Range("F15").Value2 = 1.051443266297988
$ws.Range("I15").Value2 = 1.041129907471581
$ws.Range("J15").Value2 = 1.028693225553792
$ws.Range("K15").Value2 = 1.050353928732382
$ws.Range("L15").Value2 = 1.026552895842193
$ws.Range("M15").Value2 = 1.054938534426762
$ws.Range("N15").Value2 = 1.013609764628109

# Row 16
$ws.Range("B16").Value2 = 1.02
$ws.Range("C16").Value2 = 1.022543427447125
$ws.Range("D16").Value2 = 1.047206265134796
$ws.Range("E16").Value2 = 1.023357707487512
$ws.Range("F16").Value2 = 1.051893701694854
$ws.Range("I16").Value2 = 1.041252389996923
$ws.Range("J16").Value2 = 1.028954050809598
$ws.Range("K16").Value2 = 1.050604309783654
$ws.Range("L16").Value2 = 1.026840891650843
$ws.Range("M16").Value2 = 1.055275511987177
$ws.Range("N16").Value2 = 1.013696649279117

# Row 17
$ws.Range("B17").Value2 = 1.02
$ws.Range("C17").Value2 = 1.022846407254881
$ws.Range("D17").Value2 = 1.047434799312235
$ws.Range("E17").Value2 = 1.0236120653534
$ws.Range("F17").Value2 = 1.052176515846446
$ws.Range("I17").Value2 = 1.041328884255898
$ws.Range("J17").Value2 = 1.02911767377914
$ws.Range("K17").Value2 = 1.050761267568221
$ws.Range("L17").Value2 = 1.027021654058557
$ws.Range("M17").Value2 = 1.05548693335598
$ws.Range("N17").Value2 = 1.013751146599519

# Row 18
$ws.Range("B18").Value2 = 1.02
$ws.Range("C18").Value2 = 1.023023215720863
$ws.Range("D18").Value2 = 1.047568147199764
$ws.Range("E18").Value2 = 1.023760549413171
$ws.Range("F18").Value2 = 1.052341572367024
$ws.Range("I18").Value2 = 1.041373380741467
$ws.Range("J18").Value2 = 1.029213116966615
$ws.Range("K18").Value2 = 1.050852781977905
$ws.Range("L18").Value2 = 1.027127128966363
$ws.Range("M18").Value2 = 1.055610267203743
$ws.Range("N18").Value2 = 1.013782932683507

# Row 19
$ws.Range("B19").Value2 = 1.02
$ws.Range("C19").Value2 = 1.023083517278218
$ws.Range("D19").Value2 = 1.047613623421087
$ws.Range("E19").Value2 = 1.023811199173825
$ws.Range("F19").Value2 = 1.052397868555588
$ws.Range("I19").Value2 = 1.041388532285924
$ws.Range("J19").Value2 = 1.029245661379605
$ws.Range("K19").Value2 = 1.050883979803093
$ws.Range("L19").Value2 = 1.027163099797114
$ws.Range("M19").Value2 = 1.055652323403197
$ws.Range("N19").Value2 = 1.013793770684814

# Row 20
$ws.Range("B20").Value2 = 1.02
$ws.Range("C20").Value2 = 1.022813891534643
$ws.Range("D20").Value2 = 1.047410274801611
$ws.Range("E20").Value2 = 1.023584762588006
$ws.Range("F20").Value2 = 1.052146162637606
$ws.Range("I20").Value2 = 1.041320689687824
$ws.Range("J20").Value2 = 1.029100118089507
$ws.Range("K20").Value2 = 1.05074443123818
$ws.Range("L20").Value2 = 1.027002255892033
$ws.Range("M20").Value2 = 1.055464248243667
$ws.Range("N20").Value2 = 1.013745299681829

# Row 21
$ws.Range("B21").Value2 = 1.02
$ws.Range("C21").Value2 = 1.0219383209934
$ws.Range("D21").Value2 = 1.046749734994619
$ws.Range("E21").Value2 = 1.022850026984415
$ws.Range("F21").Value2 = 1.051328973680937
$ws.Range("I21").Value2 = 1.041098703067619
$ws.Range("J21").Value2 = 1.028627000784434
$ws.Range("K21").Value2 = 1.050290320950346
$ws.Range("L21").Value2 = 1.026479801508274
$ws.Range("M21").Value2 = 1.054852982261275
$ws.Range("N21").Value2 = 1.013587701791833

# Row 22
$ws.Range("B22").Value2 = 1.02
$ws.Range("C22").Value2 = 1.021388598086351
$ws.Range("D22").Value2 = 1.046334872349068
$ws.Range("E22").Value2 = 1.022389177315232
$ws.Range("F22").Value2 = 1.050816053682144
$ws.Range("I22").Value2 = 1.040958044850485
$ws.Range("J22").Value2 = 1.028329584540765
$ws.Range("K22").Value2 = 1.050004486552596
$ws.Range("L22").Value2 = 1.02615167755326
$ws.Range("M22").Value2 = 1.054468806248915
$ws.Range("N22").Value2 = 1.013488605510493

# Row 23
$ws.Range("B23").Value2 = 1.02
$ws.Range("C23").Value2 = 1.021679941743979
$ws.Range("D23").Value2 = 1.046554756003802
$ws.Range("E23").Value2 = 1.022633376767795
$ws.Range("F23").Value2 = 1.051087878172724
$ws.Range("I23").Value2 = 1.041032713094266
$ws.Range("J23").Value2 = 1.028487245200175
$ws.Range("K23").Value2 = 1.050156042322806
$ws.Range("L23").Value2 = 1.026325587286058
$ws.Range("M23").Value2 = 1.054672450264106
$ws.Range("N23").Value2 = 1.013541138944315

# Row 24
$ws.Range("B24").Value2 = 1.02
$ws.Range("C24").Value2 = 1.022828583725153
$ws.Range("D24").Value2 = 1.047421356225006
$ws.Range("E24").Value2 = 1.023597099157561
$ws.Range("F24").Value2 = 1.052159877650414
$ws.Range("I24").Value2 = 1.041324392835391
$ws.Range("J24").Value2 = 1.029108050735701
$ws.Range("K24").Value2 = 1.050752038963723
$ws.Range("L24").Value2 = 1.027011020967411
$ws.Range("M24").Value2 = 1.055474498622303
$ws.Range("N24").Value2 = 1.013747941656277

# Row 25
$ws.Range("B25").Value2 = 1.02
$ws.Range("C25").Value2 = 1.024164864788166
$ws.Range("D25").Value2 = 1.048428863699417
$ws.Range("E25").Value2 = 1.024720201095
$ws.Range("F25").Value2 = 1.053407623428549
$ws.Range("I25").Value2 = 1.041658108286049
$ws.Range("J25").Value2 = 1.029828639633769
$ws.Range("K25").Value2 = 1.051442225894203
$ws.Range("L25").Value2 = 1.027807963065092
$ws.Range("M25").Value2 = 1.056405826035482
$ws.Range("N25").Value2 = 1.013987872699545
